$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.235.94'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '3.081.30'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''539.49'
$ws.Range("E5").Value = '  -2.78%  '
$ws.Range("D6").Value = '''135.38'
$ws.Range("E6").Value = '  -1.26%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.073.76'
$ws.Range("E8").Value = '  -0.34%  '
$ws.Range("D9").Value = '''0.496'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").Value = '''0.155'
$ws.Range("E10").Value = '  -2.58%  '
$ws.Range("D11").Value = '''6.14'
$ws.Range("E11").Value = '  -7.46%  '
$ws.Range("D12").Value = '''0.456'
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("D13").Value = '''0.0000225'
$ws.Range("E13").Value = '  +3.59%  '
$ws.Range("D14").Value = '''34.57'
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("D15").Value = '3.570.58'
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").Value = '63.164.23'
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = '3.072.46'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '''6.68'
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").Value = '''486.66'
$ws.Range("E20").Value = '  -3.40%  '
$ws.Range("D21").Value = '''13.41'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").Value = '''0.699'
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("D23").Value = '''7.16'
$ws.Range("E23").Value = '  -1.16%  '
$ws.Range("D24").Value = '''79.36'
$ws.Range("E24").Value = '  +2.04%  '
$ws.Range("D25").Value = '''12.21'
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("D27").Value = '''2.72'
$ws.Range("E27").Value = '  -1.49%  '
$ws.Range("D28").Value = '''8.22'
$ws.Range("E28").Value = '  +0.86%  '
$ws.Range("D29").Value = '''0.998'
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '''26.13'
$ws.Range("E30").Value = '  -0.44%  '
$ws.Range("D31").Value = '''1.89'
$ws.Range("E31").Value = '  -4.24%  '
$ws.Range("D32").Value = '''1.11'
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("D33").Value = '''2.39'
$ws.Range("E33").Value = '  -4.96%  '
$ws.Range("D34").Value = '''57.11'
$ws.Range("E34").Value = '  -3.25%  '
$ws.Range("D35").Value = '''5.40'
$ws.Range("E35").Value = '  +4.72%  '
$ws.Range("D36").Value = '''6.05'
$ws.Range("E36").Value = '  +2.81%  '
$ws.Range("D37").Value = '''488.36'
$ws.Range("E37").Value = '  -6.57%  '
$ws.Range("D38").Value = '3.162.02'
$ws.Range("E38").Value = '  +3.58%  '
$ws.Range("D39").Value = '''0.0397'
$ws.Range("E39").Value = '  -3.58%  '
$ws.Range("D40").Value = '''0.0799'
$ws.Range("E40").Value = '  +1.18%  '
$ws.Range("E41").Value = '  -5.97%  '
$ws.Range("D42").Value = '''8.14'
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("D43").Value = '''2.65'
$ws.Range("E43").Value = '  +1.40%  '
$ws.Range("D44").Value = '''0.254'
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("D46").Value = '0.0₃0539'
$ws.Range("E46").Value = '  +8.77%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = '''2.04'
$ws.Range("E47").Value = '  -1.06%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '''121.70'
$ws.Range("E48").Value = '  +0.46%  '
$ws.Range("D49").Value = '''24.62'
$ws.Range("E49").Value = '  +3.76%  '
$ws.Range("E50").Value = '  +2.78%  '
$ws.Range("D51").Value = '''2.35'
$ws.Range("E51").Value = '  +0.17%  '
